$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-179 down to 8-180
$ws.Rows.Item(7).Insert()

# Fill in the new row 7 with the new data record
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44515
$ws.Cells.Item(7, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 100112012
$ws.Cells.Item(7, 7).Value = "Espinaca"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 2200
$ws.Cells.Item(7, 11).Value = 400
$ws.Cells.Item(7, 12).Value = 500
$ws.Cells.Item(7, 13).Value = 450
$ws.Cells.Item(7, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(7, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(7, 16).Value = 900
$ws.Cells.Item(7, 17).Value = 0.5
$ws.Cells.Item(7, 18).Value = "Hortaliza"
